# Merge the four split runs of the "Test Kalman filter to handle noise"
# bullet on the "Next Steps" slide into a single run.
#
# Before: "Test " / "K" / "alman " / "filter to handle noise"  (4 runs)
# After : "Test Kalman filter to handle noise"                 (1 run)
#
# The surviving run keeps the rPr (lang="en-GB" dirty="0") that already
# belonged to the "K" run, matching the target OOXML.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shape = $s.Shapes.Item(2)
$para = $shape.TextFrame.TextRange.Paragraphs(5)

# Clear the runs we don't want to keep, working from the end backwards so
# that clearing one run (which removes it from the collection) doesn't
# shift the index of runs we still need to touch.
$para.Runs(4).Text = ""
$para.Runs(3).Text = ""
$para.Runs(1).Text = ""

# The only run left is the former "K" run (rPr lang="en-GB" dirty="0");
# give it the full merged text.
$para.Runs(1).Text = "Test Kalman filter to handle noise"
